$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coupling Parameters")

# Reduce the max investment capacity per year test volume
$ws.Range("B13").Value = 1200

# Add new parameter row: realistic_candidate_capacities
$ws.Range("C16").Value = "If this is true, the real capacity (less than 1000MW) of the power plants is chosen"
$ws.Range("A16").Value = "realistic_candidate_capacities"
$ws.Range("B16").Formula = "=IF(B13>=10000,FALSE,TRUE)"

$ws.Range("C18").Select()
